$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header in F1 and fill F2:F244 based on mapping from column E values
$ws.Range("F1").Value2 = "English_gruppe"

for ($r = 2; $r -le 244; $r++) {
    $eVal = $ws.Cells.Item($r, 5).Value2
    if ($eVal -eq "Villbie") {
        $ws.Cells.Item($r, 6).Value2 = "Wild bee"
    } elseif ($eVal -eq "Blomsterflue") {
        $ws.Cells.Item($r, 6).Value2 = "Hoverfly"
    } elseif ($eVal -eq "Honningbie") {
        $ws.Cells.Item($r, 6).Value2 = "Honeybee"
    } elseif ($eVal -eq "Humle") {
        $ws.Cells.Item($r, 6).Value2 = "Bumble bee"
    }
}

# Copy style from column E to column F for each row
$ws.Range("F1:F244").Style = $ws.Range("E1:E244").Style

$ws.Columns.Item(6).ColumnWidth = 18

# Remove autofilter
if ($ws.AutoFilterMode) {
    $ws.AutoFilterMode = $false
}

$wb.Save()
